# C5-PowerPoint.pptx edit:
#  1. Slide 6's table switches from the custom "Table_0" style to the
#     built-in "Medium Style 2 - Accent 1" table style.
#  2. The presentation's theme colour scheme is swapped from the
#     "Integral" palette to the standard "Office" palette.

$p = $ppt.ActivePresentation

# -- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $candidate = $slide.Shapes.Item($i)
    if ($candidate.HasTable) {
        $tableShape = $candidate
        break
    }
}
$table = $tableShape.Table
$table.ApplyStyle("{D8BF2E58-6E85-444C-B3E0-5FB523C2E5A6}")

# -- 2. Theme colour scheme: Integral -> Office ---------------------------
function HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$design = $p.Designs.Item(1)
$themeColors = $design.SlideMaster.Theme.ThemeColorScheme

# Index order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#              8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink
$officePalette = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

for ($i = 1; $i -le $officePalette.Length; $i++) {
    $themeColors.Colors($i).RGB = HexToComRgb($officePalette[$i - 1])
}
